$wb = $excel.ActiveWorkbook

# --- Data sheet: append the new weekly observation row (FRED pull) ---
$data = $wb.Worksheets.Item("Data")

$data.Cells.Item(96, 1).Value = 45133
$data.Cells.Item(96, 2).Value = 3172.244

# Match the date-cell formatting used by the rest of column A (copy format
# from the row above rather than rebuilding a numFmt from scratch).
$data.Cells.Item(95, 1).Copy()
$data.Cells.Item(96, 1).PasteSpecial(-4122)

# --- SeriesInfo sheet: refresh metadata to match the new pull ---
$info = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end / observation_end / last_updated are plain
# text cells in the source file (FRED API dump) - force text formatting
# before writing so Excel doesn't auto-convert the date-looking strings
# into date serials, then drop the formatting override so the cell is left
# with no explicit style (matching the rest of the sheet).
$info.Range("B3").NumberFormat = "@"
$info.Range("B3").Value = "2023-08-03"
$info.Range("B3").ClearFormats()

$info.Range("B4").NumberFormat = "@"
$info.Range("B4").Value = "2023-08-03"
$info.Range("B4").ClearFormats()

$info.Range("B7").NumberFormat = "@"
$info.Range("B7").Value = "2023-07-26"
$info.Range("B7").ClearFormats()

$info.Range("B14").NumberFormat = "@"
$info.Range("B14").Value = "2023-07-27 15:35:02-05"
$info.Range("B14").ClearFormats()

$info.Range("B15").Value = 77
